$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.645.36'
$ws.Range('E2').Value = '  -0.81%  '
$ws.Range('D3').Value = '2.912.76'
$ws.Range('E3').Value = '  +0.89%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '354.42'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '109.88'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.51%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.570'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.45%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  +1.37%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.07'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0883'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.44%  '
$ws.Range('E12').Value = '  +0.58%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.66'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.94%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.88'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.95%  '
$ws.Range('D15').Value = '3.376.21'
$ws.Range('E15').Value = '  +1.16%  '
$ws.Range('D16').Value = '2.907.75'
$ws.Range('E16').Value = '  +0.12%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.973'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.89%  '
$ws.Range('D18').Value = '51.696.24'
$ws.Range('E18').Value = '  -0.69%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.52'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.67%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.24'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.86'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.60%  '
$ws.Range('D22').Value = '0.0₃0978'
$ws.Range('E22').Value = '  -0.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.75'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '269.61'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.20%  '
$ws.Range('E25').Value = '  +0.77%  '
$ws.Range('E26').Value = '  +12.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '27.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.43%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.26'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +13.47%  '
$ws.Range('E30').Value = '  +13.18%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '10.54'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.09%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '38.15'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.01'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.17%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '52.13'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.92%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0438'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.53%  '
$ws.Range('E36').Value = '  +0.15%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.88'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -16.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.21'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.82%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.26'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.99'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.71'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.08%  '
$ws.Range('E42').Value = '  +1.70%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.71'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '118.63'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.85%  '
$ws.Range('E45').Value = '  -2.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.51'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.76%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.42'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.69%  '
$ws.Range('D48').Value = '2.121.30'
$ws.Range('E48').Value = '  -3.26%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.246'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -7.74%  '
$ws.Range('E50').Value = '  +4.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.04'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.19%  '
